# Insert a new weekly price record as row 245 (pushing existing rows 245-264
# down to 246-265), for "Femacal de La Calera" / Espinaca.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 245..264 down by one to make room for the new record.
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row 245 with the new weekly record.
$ws.Range("A245").Value = 3
$ws.Range("B245").Value = "Femacal de La Calera"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44578
$ws.Range("E245").Value = 5
$ws.Range("F245").Value = 100112012
$ws.Range("G245").Value = "Espinaca"
$ws.Range("H245").Value = "Sin especificar"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 170
$ws.Range("K245").Value = 4000
$ws.Range("L245").Value = 4500
$ws.Range("M245").Value = 4265
$ws.Range("N245").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O245").Value = "Provincia de Quillota"
$ws.Range("P245").Value = 1422
$ws.Range("Q245").Value = 3
$ws.Range("R245").Value = "Hortaliza"
